# LinkedIn carousel draft cleanup:
#  - drop the unused empty "Title 1" placeholder on every slide
#  - rename/renumber the two surviving textboxes (TextBox 2->1, TextBox 3->2)
#  - shrink + bold the headline textbox
#  - nudge the body textbox down and shrink it to match the new headline height

$p = $ppt.ActivePresentation

# EMU -> point helper (PowerPoint's COM surface stores Left/Top/Width/Height
# as points, internally a 32-bit float). 1 pt = 12700 EMU. Because the value
# round-trips through a 32-bit float before being re-expressed as EMU on
# save, a couple of the target points truncate one EMU short of the exact
# target; add a tiny (sub-EMU-noticeable) epsilon so the truncation lands on
# the intended EMU value.
function EmuToPt($emu) {
    $pt = [double]$emu / 12700.0
    $roundTrip = [int64]([float]$pt * 12700.0)
    if ($roundTrip -lt $emu) {
        $pt += 0.00005
    }
    return $pt
}

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)

    # Remove the empty "Title 1" placeholder shape (always shape 1 on these slides).
    $s.Shapes.Item(1).Delete()

    # Former "TextBox 2" is now the first remaining shape -> becomes "TextBox 1".
    $headline = $s.Shapes.Item(1)
    $headline.Name = "TextBox 1"
    $headline.Left = EmuToPt 457200
    $headline.Top = EmuToPt 274320
    $headline.Width = EmuToPt 8229600
    $headline.Height = EmuToPt 822960
    $headline.TextFrame.TextRange.Font.Bold = $true

    # Former "TextBox 3" is now the second remaining shape -> becomes "TextBox 2".
    $body = $s.Shapes.Item(2)
    $body.Name = "TextBox 2"
    $body.Left = EmuToPt 731520
    $body.Top = EmuToPt 4389120
    $body.Width = EmuToPt 7772400
    $body.Height = EmuToPt 1645920
}
